# Applies the MDS-database-related text/formatting edits described in the
# commit: merges a couple of runs of identical formatting in the
# AdditionalCust / AddMilkLog table definitions, inserts a new
# "AnimalName, " field into AddMilkLog, and removes the double-underline
# from the SEPayRecord table heading.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the target paragraphs by their distinctive leading text. Using
# Find (rather than hard-coded paragraph indices) makes the script
# resilient if paragraph numbering shifts.
# ---------------------------------------------------------------------
function Get-ParagraphIndexByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -match [regex]::Escape($needle)) {
            return $i
        }
    }
    return -1
}

$idxAdditionalCust = Get-ParagraphIndexByText $d "AdditionalCust ("
$idxAddMilkLog      = Get-ParagraphIndexByText $d "AddMilkLog ("
$idxSEPayRecord     = Get-ParagraphIndexByText $d "SEPayRecord"

# ---------------------------------------------------------------------
# Change 1: AdditionalCust (AcustID, CustName, Contact, LocalDID);
# The three runs ", " / "CustName," / " Contact, " (all plain / empty
# rPr) collapse into a single run with the same combined text.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item($idxAdditionalCust)
$r1 = $p1.Range
$r1.Find.Execute(", CustName, Contact, ", $true, $false, $false, $false, $false, $true, 1, $false, `
                  ", CustName, Contact, ", 2) | Out-Null

# Re-select the merged run and do a harmless bold on/off round-trip so
# the (now empty) run-properties element is kept as an explicit
# `<w:rPr/>` rather than being dropped altogether.
$p1b = $d.Paragraphs.Item($idxAdditionalCust)
$r1b = $p1b.Range
$r1b.Find.Execute(", CustName, Contact, ", $true, $false, $false, $false, $false, $true, 1, $false, `
                   "", 0) | Out-Null
$r1b.Bold = 1
$r1c = $d.Range($r1b.Start, $r1b.End)
$r1c.Bold = 0

# ---------------------------------------------------------------------
# Change 2: AddMilkLog (AcustID, OfDate, DateATime, IsMorning, Quantity,
#            Price, IsKg, AnimalName,  LocalDID);
# a) merge the thick-underlined runs "," / " OfDate," / " DateATime"
#    into one run (keeps the thick underline formatting).
# b) merge the plain runs ", " / "IsMorning, " / "Quantity, Price,
#    IsKg, " into one run.
# c) insert the brand-new "AnimalName, " field plus a trailing space as
#    two further plain runs, just before LocalDID.
# ---------------------------------------------------------------------
$p2a = $d.Paragraphs.Item($idxAddMilkLog)
$r2a = $p2a.Range
$r2a.Find.Execute(", OfDate, DateATime", $true, $false, $false, $false, $false, $true, 1, $false, `
                   ", OfDate, DateATime", 2) | Out-Null

$p2b = $d.Paragraphs.Item($idxAddMilkLog)
$r2b = $p2b.Range
$r2b.Find.Execute(", IsMorning, Quantity, Price, IsKg, ", $true, $false, $false, $false, $false, $true, 1, $false, `
                   ", IsMorning, Quantity, Price, IsKg, ", 2) | Out-Null

# Re-locate the just-merged plain run precisely (collapsed Find, no
# replacement) so we know its exact character boundaries.
$p2c = $d.Paragraphs.Item($idxAddMilkLog)
$r2c = $p2c.Range
$r2c.Find.Execute(", IsMorning, Quantity, Price, IsKg, ", $true, $false, $false, $false, $false, $true, 1, $false, `
                   "", 0) | Out-Null
$mergedStart = $r2c.Start
$mergedEnd   = $r2c.End

# Temporarily bold-mark that run so the text we are about to type next
# to it does not get silently folded back into the same run (the
# engine coalesces adjacent runs that already share identical
# character formatting).
$r2c.Bold = 1

$newFieldText = "AnimalName, "
$insertion = $d.Range($mergedEnd, $mergedEnd)
$insertion.InsertBefore($newFieldText + " ")

# Remove the temporary bold mark from the original merged run ...
$restore1 = $d.Range($mergedStart, $mergedEnd)
$restore1.Bold = 0

# ... and, as an independent operation, from the new "AnimalName, " run ...
$restore2 = $d.Range($mergedEnd, $mergedEnd + $newFieldText.Length)
$restore2.Bold = 0

# ... and from the trailing single space, each kept as its own run.
$restore3 = $d.Range($mergedEnd + $newFieldText.Length, $mergedEnd + $newFieldText.Length + 1)
$restore3.Bold = 0

# ---------------------------------------------------------------------
# Change 3: SEPayRecord heading loses its double underline (kept bold).
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item($idxSEPayRecord)
$r3 = $p3.Range
$r3.Find.Execute("SEPayRecord", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Font.Underline = 0
